$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.285.78'
$ws.Range("E2").Value = '  -1.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.878.68'
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.43'
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4844'
$ws.Range("E7").Value = '  -0.86%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2875'
$ws.Range("E8").Value = '  -3.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06588'
$ws.Range("E9").Value = '  -2.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.881.86'
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.73'
$ws.Range("E11").Value = '  -1.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07287'
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.211'
$ws.Range("E13").Value = '  +1.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.22'
$ws.Range("E14").Value = '  -2.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6556'
$ws.Range("E15").Value = '  -2.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.262.70'
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9993'
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007718'
$ws.Range("E19").Value = '  -2.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.128.05'
$ws.Range("E20").Value = '  +0.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.312'
$ws.Range("E21").Value = '  +6.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '195.46'
$ws.Range("E23").Value = '  -6.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.128'
$ws.Range("E24").Value = '  -1.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.306'
$ws.Range("E25").Value = '  -3.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.45'
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.06'
$ws.Range("E27").Value = '  -4.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.916'
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.440'
$ws.Range("E29").Value = '  +1.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.268'
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09139'
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.066'
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05119'
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7159'
$ws.Range("E34").Value = '  -4.48%  '
$ws.Range("E35").Value = '  -1.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.711'
$ws.Range("E36").Value = '  +1.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01802'
$ws.Range("E37").Value = '  -1.92%  '
$ws.Range("E38").Value = '  -1.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9187'
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.039'
$ws.Range("E40").Value = '  -2.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '106.33'
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4283'
$ws.Range("E42").Value = '  -4.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.793'
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9988'
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.24'
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.399'
$ws.Range("E46").Value = '  -5.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1321'
$ws.Range("E47").Value = '  -3.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.134'
$ws.Range("E48").Value = '  +2.90%  '
$ws.Range("E49").Value = '  -2.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05755'
$ws.Range("E50").Value = '  -2.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3834'
$ws.Range("E51").Value = '  -6.14%  '
